# Funcionarios.xlsx - add newly-submitted employee records.
# Mirrors the "Codigo de inserir a partir de excel refatorado" commit:
# eight new rows of employee data (two duplicated "import batches" of the
# same four people) are appended below the existing table, plus a couple
# of formatting touches (an underlined placeholder cell for the next
# manual entry, and a style nudge on the last new data row) that were
# left behind by the editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-EmployeeRow($r, $re, $name, $surname, $phone, $department, $role) {
  $ws.Cells.Item($r, 1).Value = $re
  $ws.Cells.Item($r, 2).Value = $name
  $ws.Cells.Item($r, 3).Value = $surname
  $ws.Cells.Item($r, 4).Value = $phone
  $ws.Cells.Item($r, 5).Value = $department
  $ws.Cells.Item($r, 6).Value = $role
}

# First import batch (rows 7-10). Row 10 is written before row 9 so the
# new shared-string entries land in the same order as the source workbook
# (Agua, Arroz, Teste, Malaquias).
Set-EmployeeRow 7  6252 "Jeferson" "Irineu" 1331241241     "Banana" "Nabo"
Set-EmployeeRow 8  3131 "Gabriel"  "Jesus"  3135515135     "Banana" "Nabo"
Set-EmployeeRow 10 455  "Marcos"   "Polo"   33552334554657 "Agua"   "Arroz"
Set-EmployeeRow 9  313  "Igor"     "Banaa"  3113513513     "Teste"  "Malaquias"

# Second import batch (rows 11-14) - same four people, new phone numbers.
Set-EmployeeRow 11 6546  "Jeferson" "Irineu" 562434414   "Banana" "Nabo"
Set-EmployeeRow 12 9078  "Gabriel"  "Jesus"  457984664   "Banana" "Nabo"
Set-EmployeeRow 13 875   "Igor"     "Banaa"  467425644   "Teste"  "Malaquias"
Set-EmployeeRow 14 78456 "Marcos"   "Polo"   4675434355  "Agua"   "Arroz"

# Placeholder cell for the next manual entry, underlined, and selected -
# this is where the cursor was left when the workbook was saved.
$ws.Range("D16").Font.Underline = $true

# Stray formatting touch on row 7 that widened the sheet's used range out
# to column I (picked up by the row's `spans` / the sheet `dimension`).
$ws.Cells.Item(7, 9).Font.Underline = $false

[void]$ws.Range("D16").Select()
